$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

$data = @{
    2 = @(-17.01697146940842,0.05797484415485412,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    3 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.653549935633869,-17.01697146940842,-17.01697146940842)
    4 = @(-17.01697146940842,0.04073474809928782,0.4694000684824108,-17.01697146940842,3.832144189915268,-17.01697146940842,1.467745238641354,-17.01697146940842,2.622265004718967,-17.01697146940842)
    5 = @(-17.01697146940842,0.6710124083698032,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    6 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    7 = @(3.172622680287332,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    8 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,1.612516782618983,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    9 = @(3.457202840629995,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    10 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.212321466696309,-17.01697146940842,2.126234844942029)
    11 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,2.114224848101845,-17.01697146940842,4.321917760592424,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.702968765737426)
    12 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    13 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,1.792988160339687,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.768611903081694,2.391882748661596)
    14 = @(-17.01697146940842,-17.01697146940842,1.092966281417879,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.751720092918688)
    15 = @(-17.01697146940842,-17.01697146940842,-0.09094192704764065,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842)
    16 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.980199740183321,-17.01697146940842)
    17 = @(-17.01697146940842,0.5925319424889315,0.2673194723958939,-17.01697146940842,-17.01697146940842,-17.01697146940842,2.31810697833861,0.7416746732609414,1.74693139229924,-17.01697146940842)
    18 = @(-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.930154613843023,0.7257129617975879,1.647798702022258,-17.01697146940842)
    19 = @(-17.01697146940842,-17.01697146940842,2.923458275155649,-17.01697146940842,-17.01697146940842,-17.01697146940842,1.287954655128825,1.521857807499339,-17.01697146940842,-17.01697146940842)
    20 = @(-17.01697146940842,3.07817440541384,2.75555736108923,-17.01697146940842,2.525374632788618,-17.01697146940842,1.502637910778166,3.059923701183464,-17.01697146940842,1.911653887268988)
    21 = @(-17.01697146940842,2.674623828457146,-17.01697146940842,3.193320555892567,-17.01697146940842,-17.01697146940842,1.660311163443162,-17.01697146940842,-17.01697146940842,-17.01697146940842)
}

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $rowValues[$i]
    }
}

Write-Host "applied supplemental-figure PSSM values"
